$d = $word.ActiveDocument

# ---------------------------------------------------------------------
# 1) Insert a new bold paragraph "Действующее лицо: гражданин." right
#    after the first paragraph ("Оплата проживания."), before the
#    "Общее предусловие..." paragraph.
# ---------------------------------------------------------------------
$secondPara = $d.Paragraphs.Item(2)
$insertionPoint = $d.Range($secondPara.Range.Start, $secondPara.Range.Start)
$insertionPoint.InsertBefore("Действующее лицо: гражданин." + [char]13)

$newPara = $d.Paragraphs.Item(2)
# Bold the whole new paragraph (text + paragraph mark) first ...
$newPara.Range.Font.Bold = 1
# ... then remove the bold from the visible text only, leaving the
# paragraph-mark formatting (w:pPr/w:rPr) bold, matching a heading-style
# paragraph whose run itself carries no direct bold formatting.
$newParaTextOnly = $d.Range($newPara.Range.Start, $newPara.Range.End - 1)
$newParaTextOnly.Font.Bold = 0

# ---------------------------------------------------------------------
# 2) Collapse the three runs "Транзакция " / "прошла неудачно" / "."
#    into a single run reading "Транзакция прошла неудачно."
# ---------------------------------------------------------------------
$d.Content.Find.Execute("Транзакция прошла неудачно.", $true, $false, $false, $false, $false, $true, 1, $false, "Транзакция прошла неудачно.", 2)

# ---------------------------------------------------------------------
# 3) Collapse the three runs "Заполнены " / "не " / "все поля. " into a
#    single run reading "Заполнены не все поля. "
# ---------------------------------------------------------------------
$d.Content.Find.Execute("Заполнены не все поля. ", $true, $false, $false, $false, $false, $true, 1, $false, "Заполнены не все поля. ", 2)
